$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered alignment) from H1
# onto the two new header cells so they match the existing header row style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for columns I (I0) and J (IF), rows 2-20
$data = @{
    2  = @(1, 3)
    3  = @(1, 5)
    4  = @(1, 5)
    5  = @(1, 7)
    6  = @(1, 4)
    7  = @(1, 6)
    8  = @(1, 2)
    9  = @(1, 6)
    10 = @(1, 6)
    11 = @(1, 6)
    12 = @(1, 3)
    13 = @(5, 7)
    14 = @(1, 5)
    15 = @(1, 6)
    16 = @(4, 5)
    17 = @(1, 5)
    18 = @(1, 4)
    19 = @(1, 3)
    20 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
